# Auto-save via app Streamlit
# Applies:
#  1. C2: "nan" -> empty
#  2. Insert a new row at row 3 (Barbara Pieper booking), shifting all
#     subsequent rows (old 3..52) down to (4..53).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear the "nan" placeholder phone number on row 2.
$ws.Range("C2").Value = ""

# 2. Insert a new row before row 3; existing rows 3-52 shift to 4-53.
$ws.Rows.Item(3).Insert()

# 3. Populate the newly inserted row 3 with the Barbara Pieper booking.
$ws.Range("A3").Value = "Barbara Pieper"
$ws.Range("B3").Value = "Booking"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = 45548
$ws.Range("E3").Value = 45881
$ws.Range("F3").Value = 333
$ws.Range("G3").Value = 274
$ws.Range("H3").Value = 223.47
$ws.Range("I3").Value = 50.53
$ws.Range("J3").Value = 18.44
$ws.Range("K3").Value = 2024
$ws.Range("L3").Value = 9
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = ""
